# Auto-generated edit script for 广州-漫展信息.xlsx
# Updates 'F' (想去人数) and some 'G' (price) columns to the new crawl snapshot values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 13798
$ws.Range("F4").Value = 41
$ws.Range("F6").Value = 805
$ws.Range("G6").Value = 70
$ws.Range("F7").Value = 2199
$ws.Range("F8").Value = 208
$ws.Range("F9").Value = 138
$ws.Range("F10").Value = 119
$ws.Range("F11").Value = 254
$ws.Range("G11").Value = 55
$ws.Range("F13").Value = 621
$ws.Range("F14").Value = 470
$ws.Range("F15").Value = 541
$ws.Range("F16").Value = 351
$ws.Range("F17").Value = 38
$ws.Range("F18").Value = 319
$ws.Range("F19").Value = 907
$ws.Range("F20").Value = 170
$ws.Range("F22").Value = 49
$ws.Range("F23").Value = 4
$ws.Range("F25").Value = 119
$ws.Range("F26").Value = 47

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 60
$ws.Range("F6").Value = 140
$ws.Range("F8").Value = 2148
$ws.Range("F11").Value = 2
$ws.Range("F13").Value = 90
$ws.Range("F15").Value = 1918

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 210
$ws.Range("F4").Value = 141

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 13798
$ws.Range("F5").Value = 41
$ws.Range("F7").Value = 805
$ws.Range("G7").Value = 70
$ws.Range("F8").Value = 60
$ws.Range("F10").Value = 2199
$ws.Range("F11").Value = 210
$ws.Range("F12").Value = 208
$ws.Range("F13").Value = 138
$ws.Range("F14").Value = 119
$ws.Range("F15").Value = 254
$ws.Range("G15").Value = 55
$ws.Range("F19").Value = 140
$ws.Range("F20").Value = 141
$ws.Range("F21").Value = 621
$ws.Range("F22").Value = 470
$ws.Range("F23").Value = 541
$ws.Range("F24").Value = 351
$ws.Range("F25").Value = 38
$ws.Range("F26").Value = 319
$ws.Range("F27").Value = 907
$ws.Range("F29").Value = 2148
$ws.Range("F32").Value = 2
$ws.Range("F34").Value = 170
$ws.Range("F36").Value = 49
$ws.Range("F37").Value = 4
$ws.Range("F38").Value = 90
$ws.Range("F41").Value = 119
$ws.Range("F42").Value = 47
$ws.Range("F43").Value = 1918
